$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "有效管理并复现关键帧（高）" paragraph: the leading _GoBack bookmark is
#    dropped from here. We achieve this by re-anchoring the (uniquely named)
#    _GoBack bookmark elsewhere in step 4 below - Bookmarks.Add() re-targets
#    the existing bookmark instead of creating a duplicate, so it disappears
#    from its old spot automatically.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 2) UI paragraph: "队伍颜色的UI：放入主界面 + 主菜单的Options（高）"
#    becomes "UI：队伍颜色，帮助等（高）".
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("队伍颜色的UI：放入主界面 + 主菜单的Options", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = "UI：队伍颜色，帮助等"

# ---------------------------------------------------------------------------
# 3) 运动系统 paragraph: "...阻尼等（中）" becomes "...阻尼等（高）".
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("另外考虑空气与水的阻尼等（中）", $true, $false, $false, $false, $false, $true, 1, $false, "另外考虑空气与水的阻尼等（高）", 2)

# ---------------------------------------------------------------------------
# 4) Swap the "雷雨等天气效果（低）" / "为MiniMap添加边框（低）" paragraphs:
#    the MiniMap item moves up (priority becomes 中) and picks up the
#    relocated _GoBack bookmark just before its closing "）"; the weather
#    item moves down, unchanged apart from position.
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("为MiniMap添加边框（低）", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = "雷雨等天气效果（低）"

$r = $d.Content
$null = $r.Find.Execute("雷雨等天气效果（低）", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = "为MiniMap添加边框（中）"

$r = $d.Content
$null = $r.Find.Execute("为MiniMap添加边框（中）", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPos = $r.Start + 14
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
